$p = $ppt.ActivePresentation

# --- 1. Update the notes master date field (2/13/2023 -> 2/16/2023) ---
$nm = $p.NotesMaster
$found = $false
foreach ($shp in $nm.Shapes) {
    if ($shp.HasTextFrame) {
        $tf = $shp.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text -eq "2/13/2023") {
                $tr.Text = "2/16/2023"
                $found = $true
            }
        }
    }
}

# --- 2. Add a vertical straight-line separator connector to slide 7 ---
$s = $p.Slides.Item(7)

$msoConnectorStraight = 1
$beginX = 3218178 / 12700.0
$beginY = 1349375 / 12700.0
$endX = $beginX
$endY = (1349375 + 2918223) / 12700.0

$conn = $s.Shapes.AddConnector($msoConnectorStraight, $beginX, $beginY, $endX, $endY)
$conn.Name = "Straight Connector 4"
$conn.Line.ForeColor.ObjectThemeColor = 13  # msoThemeColorText1 (tx1)
